$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("yyymmdd")

# The header in G1 changed from "上繳日" to "繳息迄日"
$ws.Range("G1").Value = "繳息迄日"

# Reflect the saved selection state (user had selected G1 before saving)
$ws.Range("G1").Select()
